$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 34 (this shifts existing rows 34:75 down to 35:76,
# and the sheet dimension grows from A1:T75 to A1:T76), matching the diff.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new Mango price-report record.
# (columns A,B,C,E,F,G,H,I,J,K,L,Q,T repeat the same reference/catalog values
# used throughout this data block; D,M,N,O,P,R,S carry the new record's data)
$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = 44467
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100108
$ws.Cells.Item(34, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(34, 9).Value = 100108002
$ws.Cells.Item(34, 10).Value = "Mango"
$ws.Cells.Item(34, 11).Value = "Sin especificar"
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 240
$ws.Cells.Item(34, 14).Value = 8000
$ws.Cells.Item(34, 15).Value = 8000
$ws.Cells.Item(34, 16).Value = 8000
$ws.Cells.Item(34, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(34, 18).Value = "Brasil"
$ws.Cells.Item(34, 19).Value = 2000
$ws.Cells.Item(34, 20).Value = 4
